$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Match summary header lines ---
$ws.Range("A2").Value = "England: 18/0"
$ws.Range("A3").Value = "England win by 10 wicket(s)!"

# --- England Bowling table: row 21 player renamed ---
$ws.Range("A21").Value = "Chris Jordan"

# --- England Batting table: row 26 (Jos Buttler) stats updated ---
$ws.Range("C26").Value = 18
$ws.Range("D26").Value = 3
$ws.Range("F26").Value = 3
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "600.00"
$ws.Range("G26").Style = "Normal"

# --- England Batting table: row 27 (Jonny Bairstow) stats updated ---
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 0
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "0.00"
$ws.Range("G27").Style = "Normal"

# --- Australia Bowling table: row 40 becomes "Josh Hazlewood" with updated
#     figures, and the old row 41 (which held that same player's stats) is
#     removed entirely, shifting the rest of the sheet up by one row.
$ws.Range("A40").Value = "Josh Hazlewood"
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "0.3"
$ws.Range("B40").Style = "Normal"
$ws.Range("D40").Value = 18
$ws.Range("H40").Value = 3

$ws.Rows("41").Delete()
